# "Updated to do list"
#
# On the "To Do- FY15 Release" sheet, insert a new task row above the old
# row 16 (pushing it, and everything below it, down by one):
#   Status=Done | What=Fix bug in performance ratio | Who=Janine | Priority=A
# and leave that sheet active with A17 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To Do- FY15 Release")

# Make this sheet the active tab (it was "SAM Variable Changes" before).
$ws.Activate()

# Insert a new row above the old row 16, shifting it (and every row below)
# down by one. Excel copies the formatting of the row above into the new
# row, which is what the target workbook shows (style ids 32/47/47/47/47,
# same as row 15).
$ws.Rows.Item(16).Insert()

# Fill in the new task.
$ws.Range("A16").Value = "Done"
$ws.Range("B16").Value = "Fix bug in performance ratio"
$ws.Range("C16").Value = "Janine"
$ws.Range("E16").Value = "A"

# The conditional-formatting rules below row 16 are anchored to absolute
# row numbers (A102, A104, A215, ...). A real Excel "insert row" shifts
# those along with the data; this engine's Rows.Insert doesn't touch them
# automatically, so nudge each affected rule down by one row by hand.
$fcs = $ws.Cells.FormatConditions

$snapshot = @()
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    $snapshot += , @($i, $fc.AppliesTo.Address(), $fc.Formula1)
}

foreach ($entry in $snapshot) {
    $idx = $entry[0]
    $addr = $entry[1]
    $formula = $entry[2]
    $fc = $fcs.Item($idx)

    switch ($addr) {
        '$A$1:$A$101' {
            # A1:A101 -> A1:A102 (formula still anchored at A1, unaffected)
            $fc.ModifyAppliesToRange($ws.Range("A1:A102"))
        }
        '$A$1:$A$212' {
            # A1:A212 -> A1:A213
            $fc.ModifyAppliesToRange($ws.Range("A1:A213"))
        }
        '$A$102:$A$103' {
            # A102:A103 -> A103:A104
            $fc.ModifyAppliesToRange($ws.Range("A103:A104"))
            $fc.Formula1 = $formula.Replace("A102", "A103")
        }
        '$A$215:$A$1048576' {
            # A215:A1048576 -> A216:A1048576
            $fc.ModifyAppliesToRange($ws.Range("A216:A1048576"))
            $fc.Formula1 = $formula.Replace("A215", "A216")
        }
        '$A$104:$A$212' {
            # A104:A212 (part of "A104:A212 A215:A1048576") -> A105:A213
            $fc.ModifyAppliesToRange($ws.Range("A105:A213"))
            $fc.Formula1 = $formula.Replace("A104", "A105")
        }
        default {
            # E5:E7 (and anything else) is untouched by the insert.
        }
    }
}

# Park the selection where the author left it.
$ws.Range("A17").Select()
